$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as TEXT (inline/shared string), matching
# the source file's convention of storing numeric-looking figures as text
# (t="inlineStr"). Using a leading apostrophe forces Excel to store the
# value as text; resetting the cell Style afterwards drops the transient
# "Text" number-format style that gets attached to quote-prefixed cells so
# the cell keeps the workbook's default style.
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Insert a new worksheet "2022-Q3" right after "总计", using "2022-Q2" as
#    the formatting template (same columns/header/styles), positioned before
#    it in the tab order.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, $null)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Row 2: 210009 金鹰核心资源混合 (code/name unchanged, metrics updated)
Set-TextValue $q3 2 4 "2.84"
Set-TextValue $q3 2 5 "93.42"
Set-TextValue $q3 2 6 "5.03"
Set-TextValue $q3 2 7 "0.1429"
$q3.Cells.Item(2, 8).Value = 9

# Row 3: now 001167 金鹰科技创新股票
Set-TextValue $q3 3 2 "001167"
$q3.Cells.Item(3, 3).Value = "金鹰科技创新股票"
Set-TextValue $q3 3 4 "2.66"
Set-TextValue $q3 3 5 "94.84"
Set-TextValue $q3 3 6 "5.36"
Set-TextValue $q3 3 7 "0.1426"
# H3 stays 7

# Row 4: now 162102 金鹰中小盘精选混合
Set-TextValue $q3 4 2 "162102"
$q3.Cells.Item(4, 3).Value = "金鹰中小盘精选混合"
Set-TextValue $q3 4 4 "3.17"
Set-TextValue $q3 4 5 "78.28"
Set-TextValue $q3 4 6 "4.50"
Set-TextValue $q3 4 7 "0.1426"
$q3.Cells.Item(4, 8).Value = 8

# Row 5: 005495 创金合信科技成长主题股票A (code/name unchanged, metrics updated)
Set-TextValue $q3 5 4 "1.59"
Set-TextValue $q3 5 5 "86.70"
Set-TextValue $q3 5 6 "2.81"
Set-TextValue $q3 5 7 "0.0447"
$q3.Cells.Item(5, 8).Value = 9

# Row 6: 005496 创金合信科技成长主题股票C (code/name unchanged, metrics updated)
Set-TextValue $q3 6 4 "0.62"
Set-TextValue $q3 6 5 "86.70"
Set-TextValue $q3 6 6 "2.81"
Set-TextValue $q3 6 7 "0.0174"
$q3.Cells.Item(6, 8).Value = 9

# ---------------------------------------------------------------------------
# 2. "总计" summary sheet: add the new 2022-Q3 row at the top of the data
#    (row 2) and push the rest of the quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 9 is brand new (dimension grows from D8 to D9) - clone the per-row
# number format (bold + border on column A) from row 8 before filling it in.
$total.Range("A8").Copy()
$total.Range("A9").PasteSpecial(-4122)

$rows = @(
    @(0, "2022-Q3", 5, 0.49),
    @(1, "2022-Q2", 5, 0.53),
    @(2, "2022-Q1", 8, 1.31),
    @(3, "2021-Q4", 12, 0.62),
    @(4, "2021-Q3", 6, 1.87),
    @(5, "2021-Q2", 9, 2.92),
    @(6, "2021-Q1", 24, 10.83),
    @(7, "2020-Q4", 69, 18.32)
)

$r = 2
foreach ($row in $rows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
